$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.348.36"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "3.637.48"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "578.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").Value = "3.632.49"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "4.227.08"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "3.645.39"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").Value = "68.377.13"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +26.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.99%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "690.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.422"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.22%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "0.0₃0785"
$ws.Range("E40").Value = "  +6.35%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.138"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.41%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.48%  "
$ws.Range("D43").Value = "3.225.53"
$ws.Range("E43").Value = "  +19.36%  "
$ws.Range("E44").Value = "  +13.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +31.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0421"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
